$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.254.51"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").Value = "1.607.26"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.0000"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.16"
$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3767"
$ws.Range("E7").Value = "  -0.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.51"
$ws.Range("E8").Value = "  +5.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3635"
$ws.Range("E9").Value = "  -0.21%  "

$ws.Range("E10").Value = "  +1.14%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08156"
$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.94"
$ws.Range("E13").Value = "  +1.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.598"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.407"
$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001253"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").Value = "1.605.87"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.05"
$ws.Range("E18").Value = "  +2.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06930"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.18"
$ws.Range("E20").Value = "  -0.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.549"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.94"
$ws.Range("E23").Value = "  -0.95%  "

$ws.Range("D24").Value = "23.247.21"
$ws.Range("E24").Value = "  +0.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.447"
$ws.Range("E25").Value = "  +3.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.071"
$ws.Range("E26").Value = "  +8.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.22"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.01"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.277"
$ws.Range("E29").Value = "  +0.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.53"
$ws.Range("E30").Value = "  +1.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.381"
$ws.Range("E31").Value = "  +2.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.797"
$ws.Range("E32").Value = "  -0.43%  "

$ws.Range("D33").Value = "1.780.54"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9671"
$ws.Range("E34").Value = "  +0.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07495"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("E36").Value = "  +2.30%  "

$ws.Range("E37").Value = "  +0.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2523"
$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.157"
$ws.Range("E39").Value = "  -1.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08806"
$ws.Range("E40").Value = "  -0.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.424"
$ws.Range("E41").Value = "  +4.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7118"
$ws.Range("E42").Value = "  +1.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.56"
$ws.Range("E43").Value = "  +1.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.74"
$ws.Range("E44").Value = "  +3.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6564"
$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.335"
$ws.Range("E46").Value = "  +1.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.010"

$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07950"
$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.211"
$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("E51").Value = "  -3.23%  "
